$d = $word.ActiveDocument

# Table 2: "Materias mas reprobadas" (ASIGNATURA / NO. DE ALUMNOS REPROBADOS / PORCENTAJE DE REPROBADOS)
$t1 = $d.Tables.Item(2)
$t1.Cell(3, 1).Range.Text = "ÉTICA"
$t1.Cell(3, 2).Range.Text = "12"
$t1.Cell(3, 3).Range.Text = "42.9%"

$t1.Cell(4, 1).Range.Text = "GEOMETRÍA ANALÍTICA"
$t1.Cell(4, 2).Range.Text = "9"
$t1.Cell(4, 3).Range.Text = "32.1%"

$t1.Cell(5, 1).Range.Text = "APLICA LA METODOLOGÍA ESPIRAL CON PROGRAMACIÓN ORIENTADA A OBJETOS"

$t1.Cell(6, 1).Range.Text = "INGLÉS III"
$t1.Cell(6, 2).Range.Text = "4"
$t1.Cell(6, 3).Range.Text = "14.3%"

# Table 3: "Asesorias" (ASIGNATURA DE ASESORÍA / NO. DE ALUMNOS CON ASESORÍA / NOMBRE DEL ASESOR Y/O ALUMNO MONITOR)
$t2 = $d.Tables.Item(3)
$t2.Cell(3, 1).Range.Text = "ÉTICA"
$t2.Cell(3, 2).Range.Text = "12"
$t2.Cell(3, 3).Range.Text = "Delfina Hernández Mendoza"

$t2.Cell(4, 1).Range.Text = "GEOMETRÍA ANALÍTICA"
$t2.Cell(4, 2).Range.Text = "9"
$t2.Cell(4, 3).Range.Text = "Salvador Muñoz Rivadeneyra"

$t2.Cell(5, 1).Range.Text = "APLICA LA METODOLOGÍA ESPIRAL CON PROGRAMACIÓN ORIENTADA A OBJETOS"
$t2.Cell(5, 3).Range.Text = "Miguel Sánchez Sánchez"

$t2.Cell(6, 1).Range.Text = "INGLÉS III"
$t2.Cell(6, 2).Range.Text = "4"
$t2.Cell(6, 3).Range.Text = "Luis Arturo Villanueva Morales"

# Table 4: Totales (Numero de Alumnos que tienen asesoria / canalizados / no requirieron atencion)
$t3 = $d.Tables.Item(4)
$t3.Cell(1, 5).Range.Text = "6"
$t3.Cell(3, 5).Range.Text = "7"
